# "Fruta / hortaliza, semanal" - adds a new weekly price record.
# A new row is inserted at row 159 (pushing the existing rows 159-265 down
# to 160-266) and populated with the new week's data for Albahaca at
# Femacal de La Calera.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 159, shifting rows 159:265 down
# to 160:266 (dimension grows from A1:R265 to A1:R266).
$ws.Rows("159").Insert()

# Populate the newly inserted row 159 with the new observation.
$ws.Cells.Item(159, 1).Value  = 3
$ws.Cells.Item(159, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(159, 3).Value  = "Coquimbo"
$ws.Cells.Item(159, 4).Value  = 45233
$ws.Cells.Item(159, 5).Value  = 5
$ws.Cells.Item(159, 6).Value  = 100112052
$ws.Cells.Item(159, 7).Value  = "Albahaca"
$ws.Cells.Item(159, 8).Value  = "Sin especificar"
$ws.Cells.Item(159, 9).Value  = "Primera"
$ws.Cells.Item(159, 10).Value = 40
$ws.Cells.Item(159, 11).Value = 5000
$ws.Cells.Item(159, 12).Value = 5000
$ws.Cells.Item(159, 13).Value = 5000
$ws.Cells.Item(159, 14).Value = "$/docena de matas"
$ws.Cells.Item(159, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(159, 16).Value = 833
$ws.Cells.Item(159, 17).Value = 6
$ws.Cells.Item(159, 18).Value = "Hortaliza"
